$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Append the new log entry as row 27
$ws.Range("A27").Value = "Demo inplannen"
$ws.Range("B27").Value = "klantenservice@testbedrijf123.nl"
$ws.Range("C27").Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$ws.Range("D27").Value = "Intern verzoek / Actie voor medewerker"
$ws.Range("E27").Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$ws.Range("F27").Value = "2025-08-14 21:33:21"
$ws.Range("G27").Value = "Nee"
$ws.Range("H27").Value = "Ja"
$ws.Range("I27").Value = "Nee"
$ws.Range("J27").Value = "Nee"

# Extend the conditional formatting ranges to cover the new row
$ws.Range("D2:D26").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D27"))
$ws.Range("G2:G26").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G27"))
$ws.Range("H2:H26").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H27"))
$ws.Range("I2:I26").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I27"))
$ws.Range("J2:J26").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J27"))

# Update the Dashboard summary count for "Intern verzoek / Actie voor medewerker"
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 19
